$wb = $excel.ActiveWorkbook

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 3796.875  # H51: was 3901.8572
$ws.Cells.Item(51, 9).Value = 3150  # I51: was 3430.625
$ws.Cells.Item(51, 10).Value = 4090.9092  # J51: was 4191.846
$ws.Cells.Item(51, 11).Value = 3150  # K51: was 3430.625
$ws.Cells.Item(51, 12).Value = 4090.9092  # L51: was 4191.846
$ws.Cells.Item(51, 13).Value = -2666  # M51: was -2946.625
$ws.Cells.Item(51, 14).Value = -5058.9092  # N51: was -5159.846

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2013.5358  # H70: was 2056.6924
$ws.Cells.Item(70, 9).Value = 2066.353  # I70: was 2109.25
$ws.Cells.Item(70, 10).Value = 1931.909  # J70: was 1972.6
$ws.Cells.Item(70, 11).Value = 6199.059  # K70: was 6327.75
$ws.Cells.Item(70, 12).Value = 5795.727000000001  # L70: was 5917.799999999999
$ws.Cells.Item(70, 13).Value = -5929.059  # M70: was -6057.75
$ws.Cells.Item(70, 14).Value = -6335.727000000001  # N70: was -6457.799999999999

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(73, 8).Value = 2013.5358  # H73: was 2056.6924
$ws.Cells.Item(73, 9).Value = 2066.353  # I73: was 2109.25
$ws.Cells.Item(73, 10).Value = 1931.909  # J73: was 1972.6
$ws.Cells.Item(73, 11).Value = 6199.059  # K73: was 6327.75
$ws.Cells.Item(73, 12).Value = 5795.727000000001  # L73: was 5917.799999999999
$ws.Cells.Item(73, 13).Value = -5263.059  # M73: was -5391.75
$ws.Cells.Item(73, 14).Value = -7667.727000000001  # N73: was -7789.799999999999

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 5204.9  # H74: was 4953.1113
$ws.Cells.Item(74, 9).Value = 4115.222  # I74: was 3854.25
$ws.Cells.Item(74, 10).Value = 6096.4546  # J74: was 5832.2
$ws.Cells.Item(74, 11).Value = 4115.222  # K74: was 3854.25
$ws.Cells.Item(74, 12).Value = 6096.4546  # L74: was 5832.2
$ws.Cells.Item(74, 13).Value = -3179.222  # M74: was -2918.25
$ws.Cells.Item(74, 14).Value = -7968.4546  # N74: was -7704.2

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 5204.9  # H77: was 4953.1113
$ws.Cells.Item(77, 9).Value = 4115.222  # I77: was 3854.25
$ws.Cells.Item(77, 10).Value = 6096.4546  # J77: was 5832.2
$ws.Cells.Item(77, 11).Value = 20576.11  # K77: was 19271.25
$ws.Cells.Item(77, 12).Value = 30482.273  # L77: was 29161
$ws.Cells.Item(77, 13).Value = -15896.11  # M77: was -14591.25
$ws.Cells.Item(77, 14).Value = -39842.273  # N77: was -38521

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 2091.5806  # H80: was 2130.1562
$ws.Cells.Item(80, 9).Value = 1433.1765  # I80: was 1451.3158
$ws.Cells.Item(80, 10).Value = 2891.0715  # J80: was 3122.3076
$ws.Cells.Item(80, 11).Value = 4299.529500000001  # K80: was 4353.9474
$ws.Cells.Item(80, 12).Value = 8673.2145  # L80: was 9366.9228
$ws.Cells.Item(80, 13).Value = -3301.529500000001  # M80: was -3355.9474
$ws.Cells.Item(80, 14).Value = -10669.2145  # N80: was -11362.9228

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 2091.5806  # H83: was 2130.1562
$ws.Cells.Item(83, 9).Value = 1433.1765  # I83: was 1451.3158
$ws.Cells.Item(83, 10).Value = 2891.0715  # J83: was 3122.3076
$ws.Cells.Item(83, 11).Value = 12898.5885  # K83: was 13061.8422
$ws.Cells.Item(83, 12).Value = 26019.6435  # L83: was 28100.7684
$ws.Cells.Item(83, 13).Value = -7906.5885  # M83: was -8069.842200000001
$ws.Cells.Item(83, 14).Value = -36003.6435  # N83: was -38084.7684

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 57974384  # H86: was 78433240
$ws.Cells.Item(86, 9).Value = 22223670  # I86: was 27779052
$ws.Cells.Item(86, 10).Value = 125006980  # J86: was 200003300
$ws.Cells.Item(86, 11).Value = 22223670  # K86: was 27779052
$ws.Cells.Item(86, 12).Value = 125006980  # L86: was 200003300
$ws.Cells.Item(86, 13).Value = -22222547  # M86: was -27777929
$ws.Cells.Item(86, 14).Value = -125009226  # N86: was -200005546

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 100700800  # H88: was 83501336
$ws.Cells.Item(88, 9).Value = 167833330  # I88: was 100201200
$ws.Cells.Item(88, 10).Value = 1995  # J88: was 2000
$ws.Cells.Item(88, 11).Value = 167833330  # K88: was 100201200
$ws.Cells.Item(88, 12).Value = 1995  # L88: was 2000
$ws.Cells.Item(88, 13).Value = -167832924  # M88: was -100200794
$ws.Cells.Item(88, 14).Value = -2807  # N88: was -2812

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 57974384  # H89: was 78433240
$ws.Cells.Item(89, 9).Value = 22223670  # I89: was 27779052
$ws.Cells.Item(89, 10).Value = 125006980  # J89: was 200003300
$ws.Cells.Item(89, 11).Value = 111118350  # K89: was 138895260
$ws.Cells.Item(89, 12).Value = 625034900  # L89: was 1000016500
$ws.Cells.Item(89, 13).Value = -111112734  # M89: was -138889644
$ws.Cells.Item(89, 14).Value = -625046132  # N89: was -1000027732

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(91, 8).Value = 100700800  # H91: was 83501336
$ws.Cells.Item(91, 9).Value = 167833330  # I91: was 100201200
$ws.Cells.Item(91, 10).Value = 1995  # J91: was 2000
$ws.Cells.Item(91, 11).Value = 167833330  # K91: was 100201200
$ws.Cells.Item(91, 12).Value = 1995  # L91: was 2000
$ws.Cells.Item(91, 13).Value = -167831926  # M91: was -100199796
$ws.Cells.Item(91, 14).Value = -4803  # N91: was -4808

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 266223.56  # H106: was 328691.25
$ws.Cells.Item(106, 9).Value = 428958.16  # I106: was 557465.8
$ws.Cells.Item(106, 10).Value = 1779.875  # J106: was 1870.4286
$ws.Cells.Item(106, 11).Value = 428958.16  # K106: was 557465.8
$ws.Cells.Item(106, 12).Value = 1779.875  # L106: was 1870.4286
$ws.Cells.Item(106, 13).Value = -428327.16  # M106: was -556834.8
$ws.Cells.Item(106, 14).Value = -3041.875  # N106: was -3132.4286

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 500  # H127: was 475
$ws.Cells.Item(127, 9).Value = 500  # I127: was 475
$ws.Cells.Item(127, 10).Value = 0  # J127: was 0
$ws.Cells.Item(127, 11).Value = 1500  # K127: was 1425
$ws.Cells.Item(127, 12).Value = 0  # L127: was 0
$ws.Cells.Item(127, 13).Value = 3460  # M127: was 3535

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1615.9512  # H132: was 1658.3954
$ws.Cells.Item(132, 9).Value = 1393.9697  # I132: was 1401.2646
$ws.Cells.Item(132, 10).Value = 2531.625  # J132: was 2629.7778
$ws.Cells.Item(132, 11).Value = 4181.909100000001  # K132: was 4203.793799999999
$ws.Cells.Item(132, 12).Value = 7594.875  # L132: was 7889.3334
$ws.Cells.Item(132, 13).Value = -1651.909100000001  # M132: was -1673.793799999999
$ws.Cells.Item(132, 14).Value = -12654.875  # N132: was -12949.3334

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 428985.4  # H137: was 727067.75
$ws.Cells.Item(137, 9).Value = 1855.4231  # I137: was 2037.7646
$ws.Cells.Item(137, 10).Value = 1817157.9  # J137: was 4835571
$ws.Cells.Item(137, 11).Value = 5566.2693  # K137: was 6113.293799999999
$ws.Cells.Item(137, 12).Value = 5451473.699999999  # L137: was 14506713
$ws.Cells.Item(137, 13).Value = -3016.2693  # M137: was -3563.293799999999
$ws.Cells.Item(137, 14).Value = -5456573.699999999  # N137: was -14511813

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2130.2246  # H138: was 2968.7925
$ws.Cells.Item(138, 9).Value = 1516.4546  # I138: was 1752.4445
$ws.Cells.Item(138, 10).Value = 2630.3333  # J138: was 3594.3428
$ws.Cells.Item(138, 11).Value = 4549.3638  # K138: was 5257.333500000001
$ws.Cells.Item(138, 12).Value = 7890.999899999999  # L138: was 10783.0284
$ws.Cells.Item(138, 13).Value = 590.6361999999999  # M138: was -117.3335000000006
$ws.Cells.Item(138, 14).Value = -18170.9999  # N138: was -21063.0284

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 670.8148  # H141: was 931.94116
$ws.Cells.Item(141, 9).Value = 670.8148  # I141: was 927.6875
$ws.Cells.Item(141, 10).Value = 0  # J141: was 1000
$ws.Cells.Item(141, 11).Value = 2012.4444  # K141: was 2783.0625
$ws.Cells.Item(141, 12).Value = 0  # L141: was 3000
$ws.Cells.Item(141, 13).Value = 3167.5556  # M141: was 2396.9375
$ws.Cells.Item(141, 14).ClearContents()  # N141: was -13360

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 55007  # H61: was 52329.55
$ws.Cells.Item(61, 9).Value = 1408.9166  # I61: was 1412.6923
$ws.Cells.Item(61, 10).Value = 146889.42  # J61: was 146889.42
$ws.Cells.Item(61, 11).Value = 1408.9166  # K61: was 1412.6923
$ws.Cells.Item(61, 12).Value = 146889.42  # L61: was 146889.42
$ws.Cells.Item(61, 13).Value = -1196.9166  # M61: was -1200.6923
$ws.Cells.Item(61, 14).Value = -147313.42  # N61: was -147313.42

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 703.5  # H88: was 981.75
$ws.Cells.Item(88, 9).Value = 600  # I88: was 0
$ws.Cells.Item(88, 10).Value = 807  # J88: was 981.75
$ws.Cells.Item(88, 11).Value = 600  # K88: was 0
$ws.Cells.Item(88, 12).Value = 807  # L88: was 981.75
$ws.Cells.Item(88, 13).Value = -194  # M88: was None
$ws.Cells.Item(88, 14).Value = -1619  # N88: was -1793.75

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 703.5  # H91: was 981.75
$ws.Cells.Item(91, 9).Value = 600  # I91: was 0
$ws.Cells.Item(91, 10).Value = 807  # J91: was 981.75
$ws.Cells.Item(91, 11).Value = 600  # K91: was 0
$ws.Cells.Item(91, 12).Value = 807  # L91: was 981.75
$ws.Cells.Item(91, 13).Value = 804  # M91: was None
$ws.Cells.Item(91, 14).Value = -3615  # N91: was -3789.75

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 55007  # H136: was 52329.55
$ws.Cells.Item(136, 9).Value = 1408.9166  # I136: was 1412.6923
$ws.Cells.Item(136, 10).Value = 146889.42  # J136: was 146889.42
$ws.Cells.Item(136, 11).Value = 4226.7498  # K136: was 4238.0769
$ws.Cells.Item(136, 12).Value = 440668.26  # L136: was 440668.26
$ws.Cells.Item(136, 13).Value = -1676.7498  # M136: was -1688.0769
$ws.Cells.Item(136, 14).Value = -445768.26  # N136: was -445768.26

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2175.5715  # H86: was 1545.85
$ws.Cells.Item(86, 9).Value = 2666.6667  # I86: was 1495.3889
$ws.Cells.Item(86, 10).Value = 1807.25  # J86: was 2000
$ws.Cells.Item(86, 11).Value = 2666.6667  # K86: was 1495.3889
$ws.Cells.Item(86, 12).Value = 1807.25  # L86: was 2000
$ws.Cells.Item(86, 13).Value = -1543.6667  # M86: was -372.3888999999999
$ws.Cells.Item(86, 14).Value = -4053.25  # N86: was -4246

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2175.5715  # H89: was 1545.85
$ws.Cells.Item(89, 9).Value = 2666.6667  # I89: was 1495.3889
$ws.Cells.Item(89, 10).Value = 1807.25  # J89: was 2000
$ws.Cells.Item(89, 11).Value = 13333.3335  # K89: was 7476.9445
$ws.Cells.Item(89, 12).Value = 9036.25  # L89: was 10000
$ws.Cells.Item(89, 13).Value = -7717.333500000001  # M89: was -1860.9445
$ws.Cells.Item(89, 14).Value = -20268.25  # N89: was -21232

# CRP row 19
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 894800.1  # H19: was 894797.75
$ws.Cells.Item(19, 9).Value = 1133396  # I19: was 1062559
$ws.Cells.Item(19, 10).Value = 65.5  # J19: was 70.666664
$ws.Cells.Item(19, 11).Value = 1133396  # K19: was 1062559
$ws.Cells.Item(19, 12).Value = 65.5  # L19: was 70.666664
$ws.Cells.Item(19, 13).Value = -1133226  # M19: was -1062389
$ws.Cells.Item(19, 14).Value = -405.5  # N19: was -410.666664

# CRP row 24
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(24, 8).Value = 894800.1  # H24: was 894797.75
$ws.Cells.Item(24, 9).Value = 1133396  # I24: was 1062559
$ws.Cells.Item(24, 10).Value = 65.5  # J24: was 70.666664
$ws.Cells.Item(24, 11).Value = 1133396  # K24: was 1062559
$ws.Cells.Item(24, 12).Value = 65.5  # L24: was 70.666664
$ws.Cells.Item(24, 13).Value = -1133226  # M24: was -1062389
$ws.Cells.Item(24, 14).Value = -405.5  # N24: was -410.666664

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3398.05  # H31: was 3643.1333
$ws.Cells.Item(31, 9).Value = 2155.3333  # I31: was 2504.2727
$ws.Cells.Item(31, 10).Value = 5262.125  # J31: was 6775
$ws.Cells.Item(31, 11).Value = 2155.3333  # K31: was 2504.2727
$ws.Cells.Item(31, 12).Value = 5262.125  # L31: was 6775
$ws.Cells.Item(31, 13).Value = -1860.3333  # M31: was -2209.2727
$ws.Cells.Item(31, 14).Value = -5852.125  # N31: was -7365

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3398.05  # H34: was 3643.1333
$ws.Cells.Item(34, 9).Value = 2155.3333  # I34: was 2504.2727
$ws.Cells.Item(34, 10).Value = 5262.125  # J34: was 6775
$ws.Cells.Item(34, 11).Value = 2155.3333  # K34: was 2504.2727
$ws.Cells.Item(34, 12).Value = 5262.125  # L34: was 6775
$ws.Cells.Item(34, 13).Value = -1953.3333  # M34: was -2302.2727
$ws.Cells.Item(34, 14).Value = -5666.125  # N34: was -7179

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 4499  # H62: was 4947
$ws.Cells.Item(62, 9).Value = 4499  # I62: was 4947
$ws.Cells.Item(62, 10).Value = 0  # J62: was 0
$ws.Cells.Item(62, 11).Value = 4499  # K62: was 4947
$ws.Cells.Item(62, 12).Value = 0  # L62: was 0
$ws.Cells.Item(62, 13).Value = -3875  # M62: was -4323

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 4499  # H65: was 4947
$ws.Cells.Item(65, 9).Value = 4499  # I65: was 4947
$ws.Cells.Item(65, 10).Value = 0  # J65: was 0
$ws.Cells.Item(65, 11).Value = 22495  # K65: was 24735
$ws.Cells.Item(65, 12).Value = 0  # L65: was 0
$ws.Cells.Item(65, 13).Value = -19375  # M65: was -21615

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 52254.75  # H134: was 252224.75
$ws.Cells.Item(134, 9).Value = 2231  # I134: was 2400
$ws.Cells.Item(134, 10).Value = 252349.75  # J134: was 335499.66
$ws.Cells.Item(134, 11).Value = 6693  # K134: was 7200
$ws.Cells.Item(134, 12).Value = 757049.25  # L134: was 1006498.98
$ws.Cells.Item(134, 13).Value = -4158  # M134: was -4665
$ws.Cells.Item(134, 14).Value = -762119.25  # N134: was -1011568.98

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 5030.143  # H137: was 6077.15
$ws.Cells.Item(137, 9).Value = 4712.5  # I137: was 7000
$ws.Cells.Item(137, 10).Value = 5157.2  # J137: was 5846.4375
$ws.Cells.Item(137, 11).Value = 14137.5  # K137: was 21000
$ws.Cells.Item(137, 12).Value = 15471.6  # L137: was 17539.3125
$ws.Cells.Item(137, 13).Value = -9037.5  # M137: was -15900
$ws.Cells.Item(137, 14).Value = -25671.6  # N137: was -27739.3125

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18, 8).Value = 12005  # H18: was 13502.5
$ws.Cells.Item(18, 9).Value = 12005  # I18: was 12005
$ws.Cells.Item(18, 10).Value = 0  # J18: was 15000
$ws.Cells.Item(18, 11).Value = 12005  # K18: was 12005
$ws.Cells.Item(18, 12).Value = 0  # L18: was 15000
$ws.Cells.Item(18, 13).Value = -11712  # M18: was -11712
$ws.Cells.Item(18, 14).ClearContents()  # N18: was -15586

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2809.1  # H80: was 100002320
$ws.Cells.Item(80, 9).Value = 2049.5  # I80: was 166668290
$ws.Cells.Item(80, 10).Value = 2999  # J80: was 3374
$ws.Cells.Item(80, 11).Value = 2049.5  # K80: was 166668290
$ws.Cells.Item(80, 12).Value = 2999  # L80: was 3374
$ws.Cells.Item(80, 13).Value = -1051.5  # M80: was -166667292
$ws.Cells.Item(80, 14).Value = -4995  # N80: was -5370

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2809.1  # H83: was 100002320
$ws.Cells.Item(83, 9).Value = 2049.5  # I83: was 166668290
$ws.Cells.Item(83, 10).Value = 2999  # J83: was 3374
$ws.Cells.Item(83, 11).Value = 10247.5  # K83: was 833341450
$ws.Cells.Item(83, 12).Value = 14995  # L83: was 16870
$ws.Cells.Item(83, 13).Value = -5255.5  # M83: was -833336458
$ws.Cells.Item(83, 14).Value = -24979  # N83: was -26854

# LTW row 23
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(23, 8).Value = 10  # H23: was 0
$ws.Cells.Item(23, 9).Value = 10  # I23: was 0
$ws.Cells.Item(23, 10).Value = 0  # J23: was 0
$ws.Cells.Item(23, 11).Value = 10  # K23: was 0
$ws.Cells.Item(23, 12).Value = 0  # L23: was 0
$ws.Cells.Item(23, 13).Value = 220  # M23: was None

# LTW row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(25, 8).Value = 49949.5  # H25: was 50000
$ws.Cells.Item(25, 9).Value = 49949.5  # I25: was 50000
$ws.Cells.Item(25, 10).Value = 0  # J25: was 0
$ws.Cells.Item(25, 11).Value = 49949.5  # K25: was 50000
$ws.Cells.Item(25, 12).Value = 0  # L25: was 0
$ws.Cells.Item(25, 13).Value = -49719.5  # M25: was -49770

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 4987.375  # H68: was 4260.2
$ws.Cells.Item(68, 9).Value = 4980  # I68: was 3501
$ws.Cells.Item(68, 10).Value = 4999.6665  # J68: was 4766.3335
$ws.Cells.Item(68, 11).Value = 4980  # K68: was 3501
$ws.Cells.Item(68, 12).Value = 4999.6665  # L68: was 4766.3335
$ws.Cells.Item(68, 13).Value = -4231  # M68: was -2752
$ws.Cells.Item(68, 14).Value = -6497.6665  # N68: was -6264.3335

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(71, 8).Value = 4987.375  # H71: was 4260.2
$ws.Cells.Item(71, 9).Value = 4980  # I71: was 3501
$ws.Cells.Item(71, 10).Value = 4999.6665  # J71: was 4766.3335
$ws.Cells.Item(71, 11).Value = 24900  # K71: was 17505
$ws.Cells.Item(71, 12).Value = 24998.3325  # L71: was 23831.6675
$ws.Cells.Item(71, 13).Value = -21156  # M71: was -13761
$ws.Cells.Item(71, 14).Value = -32486.3325  # N71: was -31319.6675

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1617  # H82: was 1598.8334
$ws.Cells.Item(82, 9).Value = 1617  # I82: was 1403
$ws.Cells.Item(82, 10).Value = 0  # J82: was 2186.3333
$ws.Cells.Item(82, 11).Value = 1617  # K82: was 1403
$ws.Cells.Item(82, 12).Value = 0  # L82: was 2186.3333
$ws.Cells.Item(82, 13).Value = -1256  # M82: was -1042
$ws.Cells.Item(82, 14).ClearContents()  # N82: was -2908.3333

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 1617  # H85: was 1598.8334
$ws.Cells.Item(85, 9).Value = 1617  # I85: was 1403
$ws.Cells.Item(85, 10).Value = 0  # J85: was 2186.3333
$ws.Cells.Item(85, 11).Value = 1617  # K85: was 1403
$ws.Cells.Item(85, 12).Value = 0  # L85: was 2186.3333
$ws.Cells.Item(85, 13).Value = -369  # M85: was -155
$ws.Cells.Item(85, 14).ClearContents()  # N85: was -4682.3333

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 240  # H23: was 1399.3334
$ws.Cells.Item(23, 9).Value = 10  # I23: was 0
$ws.Cells.Item(23, 10).Value = 316.66666  # J23: was 1399.3334
$ws.Cells.Item(23, 11).Value = 10  # K23: was 0
$ws.Cells.Item(23, 12).Value = 316.66666  # L23: was 1399.3334
$ws.Cells.Item(23, 13).Value = 219  # M23: was None
$ws.Cells.Item(23, 14).Value = -774.66666  # N23: was -1857.3334

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2271.3914  # H107: was 2142.4736
$ws.Cells.Item(107, 9).Value = 1571.5  # I107: was 1391.9166
$ws.Cells.Item(107, 10).Value = 4791  # J107: was 3429.1428
$ws.Cells.Item(107, 11).Value = 4714.5  # K107: was 4175.7498
$ws.Cells.Item(107, 12).Value = 14373  # L107: was 10287.4284
$ws.Cells.Item(107, 13).Value = -2794.5  # M107: was -2255.7498
$ws.Cells.Item(107, 14).Value = -18213  # N107: was -14127.4284

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2123.2693  # H122: was 2573.468
$ws.Cells.Item(122, 9).Value = 1563.5778  # I122: was 2022.575
$ws.Cells.Item(122, 10).Value = 5721.2856  # J122: was 5721.4287
$ws.Cells.Item(122, 11).Value = 4690.7334  # K122: was 6067.725
$ws.Cells.Item(122, 12).Value = 17163.8568  # L122: was 17164.2861
$ws.Cells.Item(122, 13).Value = -2240.7334  # M122: was -3617.725
$ws.Cells.Item(122, 14).Value = -22063.8568  # N122: was -22064.2861
